$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "TestData"

# Add new row 4 (TC_02 - Mandatory Field Error Test)
$ws.Range("A4").Value = "TC_02"
$ws.Range("B4").Value = "Mandatory Field Error Test"
$ws.Range("D4").Value = 68
$ws.Range("E4").Value = 100000
$ws.Range("F4").Value = 75000
$ws.Range("G4").Value = 500000
$ws.Range("H4").Value = 10
$ws.Range("I4").Value = 25
$ws.Range("J4").Value = 4000
$ws.Range("K4").Value = 500
$ws.Range("L4").Value = 20
$ws.Range("M4").Value = 75
$ws.Range("N4").Value = 8
$ws.Range("O4").Value = 5
$ws.Range("P4").Value = "Congratulations! You are exceeding your retirement goals"

# Add new header cells
$ws.Range("Q1").Value = "ErrorMsg"
$ws.Range("R1").Value = "FieldErrorMsg"

$ws.Range("Q4").Value = "Please fill out all required fields"
$ws.Range("R4").Value = "Input required"

# Page setup
$ws.PageSetup.Orientation = 1

# Match final selection in the saved file
[void]$ws.Range("R4").Select()

Write-Host "done"
